$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("A2").Value2 = 3015535
$ws.Range("B2").Value2 = 96333
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value2 = "Ovaliderad"
$ws.Range("C2").ClearFormats()
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "VU"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value2 = 220787
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value2 = "Knärot"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value2 = "Goodyera repens"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value2 = "(L.) R. Br."
$ws.Range("H2").ClearFormats()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value2 = "Siggeforasjön, Upl"
$ws.Range("P2").ClearFormats()
$ws.Range("Q2").Value2 = 620480.0184618242
$ws.Range("R2").Value2 = 6650653.458782602
$ws.Range("S2").Value2 = 10
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value2 = "Uppsala"
$ws.Range("T2").ClearFormats()
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value2 = "Uppsala"
$ws.Range("U2").ClearFormats()
$ws.Range("V2").NumberFormat = "@"
$ws.Range("V2").Value2 = "Uppland"
$ws.Range("V2").ClearFormats()
$ws.Range("W2").NumberFormat = "@"
$ws.Range("W2").Value2 = "Järlåsa"
$ws.Range("W2").ClearFormats()
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value2 = "2004-07-27"
$ws.Range("Y2").ClearFormats()
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value2 = "00:00"
$ws.Range("Z2").ClearFormats()
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value2 = "2004-07-27"
$ws.Range("AA2").ClearFormats()
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value2 = "00:00"
$ws.Range("AB2").ClearFormats()
$ws.Range("AD2").Value2 = $false
$ws.Range("AE2").Value2 = $false
$ws.Range("AF2").ClearContents()
$ws.Range("AG2").Value2 = $false
$ws.Range("AT2").ClearContents()
$ws.Range("AW2").NumberFormat = "@"
$ws.Range("AW2").Value2 = "Sandra Lindström"
$ws.Range("AW2").ClearFormats()
$ws.Range("AX2").NumberFormat = "@"
$ws.Range("AX2").Value2 = "Sandra Lindström"
$ws.Range("AX2").ClearFormats()
$ws.Range("AY2").ClearContents()

# ---- Row 3 ----
$ws.Range("A3").Value2 = 4219207
$ws.Range("B3").Value2 = 95518
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value2 = "Ovaliderad"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "LC"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value2 = 221945
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value2 = "Revlummer"
$ws.Range("F3").ClearFormats()
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value2 = "Lycopodium annotinum"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value2 = "L."
$ws.Range("H3").ClearFormats()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value2 = "Siggeforasjön, Upl"
$ws.Range("P3").ClearFormats()
$ws.Range("Q3").Value2 = 620480.0184618242
$ws.Range("R3").Value2 = 6650653.458782602
$ws.Range("S3").Value2 = 10
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value2 = "Uppsala"
$ws.Range("T3").ClearFormats()
$ws.Range("U3").NumberFormat = "@"
$ws.Range("U3").Value2 = "Uppsala"
$ws.Range("U3").ClearFormats()
$ws.Range("V3").NumberFormat = "@"
$ws.Range("V3").Value2 = "Uppland"
$ws.Range("V3").ClearFormats()
$ws.Range("W3").NumberFormat = "@"
$ws.Range("W3").Value2 = "Järlåsa"
$ws.Range("W3").ClearFormats()
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value2 = "2004-07-27"
$ws.Range("Y3").ClearFormats()
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value2 = "00:00"
$ws.Range("Z3").ClearFormats()
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value2 = "2004-07-27"
$ws.Range("AA3").ClearFormats()
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value2 = "00:00"
$ws.Range("AB3").ClearFormats()
$ws.Range("AD3").Value2 = $false
$ws.Range("AE3").Value2 = $false
$ws.Range("AG3").Value2 = $false
$ws.Range("AI3").ClearContents()
$ws.Range("AT3").ClearContents()
$ws.Range("AW3").NumberFormat = "@"
$ws.Range("AW3").Value2 = "Sandra Lindström"
$ws.Range("AW3").ClearFormats()
$ws.Range("AX3").NumberFormat = "@"
$ws.Range("AX3").Value2 = "Sandra Lindström"
$ws.Range("AX3").ClearFormats()
$ws.Range("AY3").ClearContents()

# ---- Row 4 ----
$ws.Range("A4").Value2 = 5971707
$ws.Range("B4").Value2 = 96253
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value2 = "Ovaliderad"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "LC"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value2 = 223597
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value2 = "Jungfru marie nycklar"
$ws.Range("F4").ClearFormats()
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value2 = "Dactylorhiza maculata subsp. maculata"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value2 = "Siggeforasjön, Upl"
$ws.Range("P4").ClearFormats()
$ws.Range("Q4").Value2 = 620480.0184618242
$ws.Range("R4").Value2 = 6650653.458782602
$ws.Range("S4").Value2 = 10
$ws.Range("T4").NumberFormat = "@"
$ws.Range("T4").Value2 = "Uppsala"
$ws.Range("T4").ClearFormats()
$ws.Range("U4").NumberFormat = "@"
$ws.Range("U4").Value2 = "Uppsala"
$ws.Range("U4").ClearFormats()
$ws.Range("V4").NumberFormat = "@"
$ws.Range("V4").Value2 = "Uppland"
$ws.Range("V4").ClearFormats()
$ws.Range("W4").NumberFormat = "@"
$ws.Range("W4").Value2 = "Järlåsa"
$ws.Range("W4").ClearFormats()
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value2 = "2004-07-27"
$ws.Range("Y4").ClearFormats()
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value2 = "00:00"
$ws.Range("Z4").ClearFormats()
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value2 = "2004-07-27"
$ws.Range("AA4").ClearFormats()
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value2 = "00:00"
$ws.Range("AB4").ClearFormats()
$ws.Range("AD4").Value2 = $false
$ws.Range("AE4").Value2 = $false
$ws.Range("AG4").Value2 = $false
$ws.Range("AT4").ClearContents()
$ws.Range("AW4").NumberFormat = "@"
$ws.Range("AW4").Value2 = "Sandra Lindström"
$ws.Range("AW4").ClearFormats()
$ws.Range("AX4").NumberFormat = "@"
$ws.Range("AX4").Value2 = "Sandra Lindström"
$ws.Range("AX4").ClearFormats()
$ws.Range("AY4").ClearContents()

# ---- Row 5 ----
$ws.Range("A5").Value2 = 85681612
$ws.Range("B5").Value2 = 55392
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value2 = "Ovaliderad"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "LC"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = 208257
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value2 = "Kopparödla"
$ws.Range("F5").ClearFormats()
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value2 = "Anguis fragilis"
$ws.Range("G5").ClearFormats()
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value2 = "Linnaeus, 1758"
$ws.Range("H5").ClearFormats()
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value2 = "1"
$ws.Range("I5").ClearFormats()
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value2 = "Siggeforasjön öst, Upl"
$ws.Range("P5").ClearFormats()
$ws.Range("Q5").Value2 = 620703.2845301379
$ws.Range("R5").Value2 = 6650915.221505931
$ws.Range("S5").Value2 = 1
$ws.Range("T5").NumberFormat = "@"
$ws.Range("T5").Value2 = "Uppsala"
$ws.Range("T5").ClearFormats()
$ws.Range("U5").NumberFormat = "@"
$ws.Range("U5").Value2 = "Uppsala"
$ws.Range("U5").ClearFormats()
$ws.Range("V5").NumberFormat = "@"
$ws.Range("V5").Value2 = "Uppland"
$ws.Range("V5").ClearFormats()
$ws.Range("W5").NumberFormat = "@"
$ws.Range("W5").Value2 = "Järlåsa"
$ws.Range("W5").ClearFormats()
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value2 = "2020-05-22"
$ws.Range("Y5").ClearFormats()
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value2 = "11:30"
$ws.Range("Z5").ClearFormats()
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value2 = "2020-05-22"
$ws.Range("AA5").ClearFormats()
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value2 = "11:30"
$ws.Range("AB5").ClearFormats()
$ws.Range("AD5").Value2 = $false
$ws.Range("AE5").Value2 = $false
$ws.Range("AF5").ClearContents()
$ws.Range("AG5").Value2 = $false
$ws.Range("AT5").ClearContents()
$ws.Range("AW5").NumberFormat = "@"
$ws.Range("AW5").Value2 = "Katarina Sjöholm"
$ws.Range("AW5").ClearFormats()
$ws.Range("AX5").NumberFormat = "@"
$ws.Range("AX5").Value2 = "Katarina Sjöholm"
$ws.Range("AX5").ClearFormats()
$ws.Range("AY5").ClearContents()

# ---- Row 6 ----
$ws.Range("A6").Value2 = 110602620
$ws.Range("B6").Value2 = 55395
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value2 = "Ovaliderad"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "LC"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = 208257
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value2 = "Kopparödla"
$ws.Range("F6").ClearFormats()
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value2 = "Anguis fragilis"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value2 = "Linnaeus, 1758"
$ws.Range("H6").ClearFormats()
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value2 = "1"
$ws.Range("I6").ClearFormats()
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value2 = "ex."
$ws.Range("J6").ClearFormats()
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value2 = "adult"
$ws.Range("K6").ClearFormats()
$ws.Range("N6").NumberFormat = "@"
$ws.Range("N6").Value2 = "observerad"
$ws.Range("N6").ClearFormats()
$ws.Range("P6").NumberFormat = "@"
$ws.Range("P6").Value2 = "Siggeforasjön, O om, Upl"
$ws.Range("P6").ClearFormats()
$ws.Range("Q6").Value2 = 620647.3993695766
$ws.Range("R6").Value2 = 6650890.352745522
$ws.Range("S6").Value2 = 5
$ws.Range("T6").NumberFormat = "@"
$ws.Range("T6").Value2 = "Uppsala"
$ws.Range("T6").ClearFormats()
$ws.Range("U6").NumberFormat = "@"
$ws.Range("U6").Value2 = "Uppsala"
$ws.Range("U6").ClearFormats()
$ws.Range("V6").NumberFormat = "@"
$ws.Range("V6").Value2 = "Uppland"
$ws.Range("V6").ClearFormats()
$ws.Range("W6").NumberFormat = "@"
$ws.Range("W6").Value2 = "Järlåsa"
$ws.Range("W6").ClearFormats()
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value2 = "2023-06-24"
$ws.Range("Y6").ClearFormats()
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value2 = "00:00"
$ws.Range("Z6").ClearFormats()
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value2 = "2023-06-24"
$ws.Range("AA6").ClearFormats()
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value2 = "00:00"
$ws.Range("AB6").ClearFormats()
$ws.Range("AD6").Value2 = $false
$ws.Range("AE6").Value2 = $false
$ws.Range("AG6").Value2 = $false
$ws.Range("AI6").NumberFormat = "@"
$ws.Range("AI6").Value2 = "Skogsbilväg"
$ws.Range("AI6").ClearFormats()
$ws.Range("AT6").ClearContents()
$ws.Range("AW6").NumberFormat = "@"
$ws.Range("AW6").Value2 = "Ulf Arup"
$ws.Range("AW6").ClearFormats()
$ws.Range("AX6").NumberFormat = "@"
$ws.Range("AX6").Value2 = "Ulf Arup, Stefan Ekman, Wenche Eide"
$ws.Range("AX6").ClearFormats()
$ws.Range("AY6").ClearContents()
